# "add VPC figures (#12)" -- applies the sensitivity-analysis / absorption / mass-balance
# task rows to the MeanModelSimulation sheet of WorkflowMean.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. The C7 cell used to reference the workbook's own file name
#    ("WorkflowMean.xlsx"); that reference is simply dropped.
# ---------------------------------------------------------------------------
$ws.Range("C7").ClearContents()

# ---------------------------------------------------------------------------
# 2. Make room for the old "tasks" block 3 rows further down (it moves from
#    rows 15-17 to rows 18-20) by inserting 3 blank rows after row 17.
#    Newly inserted rows inherit the row-17 look (A:s6 B:s7 D:E:F:s4), which
#    is exactly what rows 19/20 need.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# Re-create the old row 15 / 16 / 17 content at rows 18 / 19 / 20.
$ws.Range("B18").Value = "Following entries define tasks, which should be performed for all simulations. Only the entries of the first simulation are taken into account"
$ws.Range("A19").Value = "TaskdoVPC"
$ws.Range("B19").Value = "default plots will be generated"
$ws.Range("C19").Value = 1
$ws.Range("A20").Value = "TaskdoSensitivityAnalysis"
$ws.Range("B20").Value = "sensitivity analysis is performed"
$ws.Range("C20").Value = 0

# Row 18 needs the special "section footer" look (style indices 8/9) that the
# old row 15 used; copy it over explicitly (value already set above is kept).
$ws.Range("A15:F15").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Turn the old rows 15-17 into the new "Sensitivity" parameter block.
# ---------------------------------------------------------------------------
# Row 15 becomes a new section header, styled like the other section headers
# (rows 2 / 6 / 10 use style indices 4/5).
$ws.Range("A2:F2").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)
$ws.Range("B15").Value = "Sensitivity"

# Row 16: sensXls parameter (keeps its original A:s6 / B:s7 / D:E:F:s4 look)
$ws.Range("A16").Value = "sensXls"
$ws.Range("B16").Value = "xlsfilefor sensitivity Parameter definition; if it is empty, sheet is in this xlsfile"
$ws.Range("C16").ClearContents()

# Row 17: sensSheet parameter (same look as row 16)
$ws.Range("A17").Value = "sensSheet"
$ws.Range("B17").Value = "xlssheet for sensitivity Parameter definition; if empty first sheet is taken"
$ws.Range("C17").ClearContents()

# ---------------------------------------------------------------------------
# 4. Append the two new task rows (absorption plots / mass balance check)
#    plus one trailing blank row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = "TaskdoAbsorptionPlots"
$ws.Range("B21").Value = "absorption is plotted"
$ws.Range("C21").Value = 0

$ws.Range("B22").Value = "massbalance will be checked"
$ws.Range("A22").Value = "TaskcheckMassbalance"
$ws.Range("C22").Value = 0

# Row 23 stays blank (just A/B formatted like the other task rows).

# ---------------------------------------------------------------------------
# 5. Sheet view: selection moves to C7 and the sheet no longer needs to be
#    scrolled down (topLeftCell reset to default / A1).
# ---------------------------------------------------------------------------
$ws.Range("C7").Select()
